############################################################################
# Generate Report for Archive
#
# 1. The localization status text changes from "Ready for handoff" to
#    "In Translation" everywhere it appears (Overview sheet's per-language
#    status columns, plus each language sheet's own Status column).
# 2. The "Status" column is narrowed on every sheet that has one:
#       - Overview : columns E (zh-cn) and F (de-de)
#       - zh-cn     : column C
#       - de-de     : column C
############################################################################

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# Target width (narrower) expressed as Excel's stored "character" column
# width; closest value reachable through the ColumnWidth COM property
# (which snaps to 1/6-character increments) is obtained by setting 12.5.
$newColumnWidth = 12.5

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rows = $used.Rows.Count
    $cols = $used.Columns.Count

    for ($r = 1; $r -le $rows; $r++) {
        for ($c = 1; $c -le $cols; $c++) {
            $cell = $used.Cells.Item($r, $c)
            # NOTE: keep the string literal on the LEFT of -eq; Excel cells
            # holding booleans/numbers would otherwise force the literal to
            # be coerced to that type (PowerShell's left-operand-wins
            # comparison semantics) and produce false positives.
            if ($oldStatus -eq $cell.Value2) {
                $cell.Value = $newStatus
            }
        }
    }
}

# Narrow the "Status" / per-language status columns.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns("E").ColumnWidth = $newColumnWidth
$wsOverview.Columns("F").ColumnWidth = $newColumnWidth

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Columns("C").ColumnWidth = $newColumnWidth

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Columns("C").ColumnWidth = $newColumnWidth
